# Pouya Finance "Overview" income-statement update (kemina / rial.xlsx)
# The published table is a rolling 5-period window: the oldest period
# (12 ماهه منتهی به 1396/12, published 1398-04-23) is dropped, every
# remaining period shifts one column to the left (D<-E, E<-F, F<-G, G<-H),
# and the newest period (12 ماهه منتهی به 1401/12, published 1402-02-29)
# is appended in column H, together with a freshly-published figure for
# the previous period's publish date (1401-11-07 -> 1402-02-29 (9)).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 8: "دوره مالی" period-ending headers ---
$ws.Range("D8").Value = "12 ماهه منتهی به 1397/12"
$ws.Range("E8").Value = "12 ماهه منتهی به 1398/12"
$ws.Range("F8").Value = "12 ماهه منتهی به 1399/12"
$ws.Range("G8").Value = "12 ماهه منتهی به 1400/12"
$ws.Range("H8").Value = "12 ماهه منتهی به 1401/12"

# --- Row 9: "تاریخ انتشار" publish-date headers ---
$ws.Range("D9").Value = "1399-04-30 (7)"
$ws.Range("E9").Value = "1400-04-28 (8)"
$ws.Range("F9").Value = "1401-04-29 (9)"
$ws.Range("G9").Value = "1402-02-29 (9)"
$ws.Range("H9").Value = "1402-02-29"

# --- Row 11: فروش (Sales) ---
$ws.Range("D11").Value = 763875
$ws.Range("E11").Value = 1028797
$ws.Range("F11").Value = 2234267
$ws.Range("G11").Value = 3283410
$ws.Range("H11").Value = 5033548

# --- Row 12: بهای تمام شده کالای فروش رفته (Cost of goods sold) ---
$ws.Range("D12").Value = -550547
$ws.Range("E12").Value = -774825
$ws.Range("F12").Value = -1322818
$ws.Range("G12").Value = -2099151
$ws.Range("H12").Value = -4073861

# --- Row 13: سود (زیان) ناخالص (Gross profit) ---
$ws.Range("D13").Value = 213328
$ws.Range("E13").Value = 253972
$ws.Range("F13").Value = 911449
$ws.Range("G13").Value = 1184259
$ws.Range("H13").Value = 959687

# --- Row 14: هزینه های عمومی, اداری و تشکیلاتی (SG&A expenses) ---
$ws.Range("D14").Value = -61887
$ws.Range("E14").Value = -69376
$ws.Range("F14").Value = -94935
$ws.Range("G14").Value = -171526
$ws.Range("H14").Value = -342452

# --- Row 15: هزینه کاهش ارزش دریافتنی‌ها (was a literal "-" placeholder,
#     now a real number for every period) ---
$ws.Range("D15").Value = -4430
$ws.Range("E15").Value = -43539
$ws.Range("F15").Value = 0
$ws.Range("G15").Value = 0
$ws.Range("H15").Value = 0

# --- Row 16: خالص سایر درامدها (هزینه ها) ی عملیاتی ---
$ws.Range("D16").Value = 4644
$ws.Range("E16").Value = -20542
$ws.Range("F16").Value = -7520
$ws.Range("G16").Value = -78972
$ws.Range("H16").Value = -379117

# --- Row 17: سود (زیان) عملیاتی (Operating profit) ---
$ws.Range("D17").Value = 151655
$ws.Range("E17").Value = 120515
$ws.Range("F17").Value = 808994
$ws.Range("G17").Value = 933761
$ws.Range("H17").Value = 238118

# --- Row 18: هزینه های مالی (Financial expenses) ---
$ws.Range("D18").Value = -30449
$ws.Range("E18").Value = -33026
$ws.Range("F18").Value = -107342
$ws.Range("G18").Value = -59028
$ws.Range("H18").Value = -330348

# --- Row 19: خالص سایر درامدها و هزینه های غیرعملیاتی ---
$ws.Range("D19").Value = -13
$ws.Range("E19").Value = -5044
$ws.Range("F19").Value = 6289
$ws.Range("G19").Value = 6594
$ws.Range("H19").Value = 113977

# --- Row 20: سود (زیان) خالص عملیات در حال تداوم قبل از مالیات ---
$ws.Range("D20").Value = 121193
$ws.Range("E20").Value = 82445
$ws.Range("F20").Value = 707941
$ws.Range("G20").Value = 881327
$ws.Range("H20").Value = 21747

# --- Row 21: مالیات (Tax) ---
$ws.Range("D21").Value = -3754
$ws.Range("E21").Value = -47477
$ws.Range("F21").Value = -98395
$ws.Range("G21").Value = -113694
$ws.Range("H21").Value = -3914

# --- Row 22: سود (زیان) خالص عملیات در حال تداوم ---
$ws.Range("D22").Value = 117439
$ws.Range("E22").Value = 34968
$ws.Range("F22").Value = 609546
$ws.Range("G22").Value = 767633
$ws.Range("H22").Value = 17833

# --- Row 23: سود (زیان) عملیات متوقف شده پس از اثر مالیاتی ---
$ws.Range("D23").Value = 0
$ws.Range("E23").Value = 0
$ws.Range("F23").Value = 0
$ws.Range("G23").Value = 0
$ws.Range("H23").Value = 0

# --- Row 24: سود (زیان) خالص (Net profit) ---
$ws.Range("D24").Value = 117439
$ws.Range("E24").Value = 34968
$ws.Range("F24").Value = 609546
$ws.Range("G24").Value = 767633
$ws.Range("H24").Value = 17833

# --- Row 25: سود هر سهم پس از کسر مالیات (EPS after tax) ---
$ws.Range("D25").Value = 546
$ws.Range("E25").Value = 98
$ws.Range("F25").Value = 1703
$ws.Range("G25").Value = 603
$ws.Range("H25").Value = 14

# --- Row 26: سرمایه (Capital) ---
$ws.Range("D26").Value = 215211
$ws.Range("E26").Value = 358000
$ws.Range("F26").Value = 358000
$ws.Range("G26").Value = 1273000
$ws.Range("H26").Value = 1273000

# --- Row 27: سود هر سهم بر اساس آخرین سرمایه ---
$ws.Range("D27").Value = 92
$ws.Range("E27").Value = 27
$ws.Range("F27").Value = 479
$ws.Range("G27").Value = 603
$ws.Range("H27").Value = 14
